$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Modèle")
$ws.Activate()

# ---------------------------------------------------------------------------
# Formatting for the three new rows (25, 26, 27): columns A-G and J-L get the
# "no fill / vertical-center / wrap text" look (style index 2 in the source
# file); columns H-I keep the highlighted look used by row 24 (style index 13)
# ---------------------------------------------------------------------------
foreach ($r in 25..27) {
    $ws.Rows.Item($r).RowHeight = 73.75

    $rngMain1 = $ws.Range("A$r`:G$r")
    $rngMain1.VerticalAlignment = -4108  # xlCenter
    $rngMain1.WrapText = $true

    $rngMain2 = $ws.Range("J$r`:L$r")
    $rngMain2.VerticalAlignment = -4108  # xlCenter
    $rngMain2.WrapText = $true

    $ws.Range("H24:I24").Copy()
    $ws.Range("H$r`:I$r").PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = "fit_heatmap"
$ws.Range("C25").Value = "Sortir Y_pred qui prend en compte G et DBH pour faire à la suite un heatmap"
$ws.Range("E25").Value = "m_heatmap"
$ws.Range("D25").Value = "Y = softmax(alpha + beta*X_cr[n] + gamma*(X_cr[n])^2 + delta*Z_cr[n]+ epsilon*(Z_cr[n])^2)"
$ws.Range("F25").Value = "BD_esp_G"
$ws.Range("G25").Value = "BD_env_G"
$ws.Range("H25").Value = " 1992, 1993,1995,2005, 2008, 2013, 2016"
$ws.Range("I25").Value = 6
$ws.Range("B25").Value = 3
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = "8h23,10h43 ordi plante"
$ws.Range("L25").Value = "Après avoir fait tourner les chaines il a du mal et s'arrête"

# ---------------------------------------------------------------------------
# Row 26
# ---------------------------------------------------------------------------
$ws.Range("A26").Value = "fit_heatmap"
$ws.Range("C26").Value = "Sortir Y_pred qui prend en compte G et DBH pour faire à la suite un heatmap"
$ws.Range("E26").Value = "m_heatmap"
$ws.Range("D26").Value = "Y = softmax(alpha + beta*X_cr[n] + gamma*(X_cr[n])^2 + delta*Z_cr[n]+ epsilon*(Z_cr[n])^2)"
$ws.Range("F26").Value = "BD_esp_G"
$ws.Range("G26").Value = "BD_env_G"
$ws.Range("H26").Value = " 1992, 1993,1995,2005, 2008, 2013, 2016"
$ws.Range("I26").Value = 6
$ws.Range("B26").Value = 3
$ws.Range("J26").Value = 10
$ws.Range("K26").Value = "3h"
$ws.Range("L26").Value = "Foireux car seulement 10 itération qui ne permet pas. L'extraction des Y_préd prend vraiment bcp de temps"

# ---------------------------------------------------------------------------
# Row 27
# ---------------------------------------------------------------------------
$ws.Range("A27").Value = "fit_heatmap"
$ws.Range("C27").Value = "Sortir Y_pred qui prend en compte G et DBH pour faire à la suite un heatmap (réduction du nombre de Z_pred (0.1) et X_pred (5) calculer "
$ws.Range("D27").Value = "Y = softmax(alpha + beta*X_cr[n] + gamma*(X_cr[n])^2 + delta*Z_cr[n]+ epsilon*(Z_cr[n])^2)"
$ws.Range("E27").Value = "m_heatmap"
$ws.Range("F27").Value = "BD_esp_G"
$ws.Range("G27").Value = "BD_env_G"
$ws.Range("H27").Value = " 1992, 1993,1995,2005, 2008, 2013, 2016"
$ws.Range("I27").Value = 6
$ws.Range("B27").Value = 3
$ws.Range("J27").Value = 100
$ws.Range("L27").Value = "ça va beaucoup plus vite. Il converge même si R dit que non moi je trouve que une RHAT de 1.08 c'est bien."
$ws.Range("K27").Value = "6min"

# ---------------------------------------------------------------------------
# Selection / view state (best effort - matches the author's final selection)
# ---------------------------------------------------------------------------
$ws.Range("K27").Select()

Write-Host "edit applied"
